# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.110.91"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.475.84"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'576.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.32%  "
$ws.Range("D6").Value = "'146.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "2.475.19"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  +8.25%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "2.925.50"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "63.075.92"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "2.472.87"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D20").Value = "'11.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "'329.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  +9.19%  "
$ws.Range("D23").Value = "'4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'66.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'670.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.30%  "
$ws.Range("E27").Value = "  +13.75%  "
$ws.Range("D28").Value = "0.0₃0993"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").Value = "'8.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "'1.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "0.0₆0312"
$ws.Range("E45").Value = "  +11.05%  "
$ws.Range("D46").Value = "'150.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("E47").Value = "  +27.21%  "
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "'0.607"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("E51").Value = "  -0.03%  "
